# fix: alterar python version para 3.11.5
# Update the absenteeism data rows (2-11) with new values per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=10559; B="Marina Vieira";        C="Operações";  D="Consulta médica";    E=3; F=45088; G=6435.95},
    @{Row=3;  A=61634; B="Sr. Nathan Pereira";    C="Marketing";  D="Doença";             E=5; F=45093; G=7331.06},
    @{Row=4;  A=95656; B="Nicolas Caldeira";      C="TI";         D="Outros";             E=4; F=45106; G=7591.4},
    @{Row=5;  A=64875; B="Amanda Vieira";         C="Marketing";  D="Viagem de negócios"; E=4; F=45097; G=11421.36},
    @{Row=6;  A=65373; B="André Silva";           C="Jurídico";   D="Consulta médica";    E=5; F=45090; G=10068.43},
    @{Row=7;  A=27185; B="Thales Castro";         C="Jurídico";   D="Doença";             E=3; F=45099; G=9060.66},
    @{Row=8;  A=85424; B="Ana Lívia Monteiro";    C="Operações";  D="Consulta médica";    E=8; F=45094; G=7306.7},
    @{Row=9;  A=82436; B="Caroline Moreira";      C="Jurídico";   D="Consulta médica";    E=7; F=45090; G=11444.94},
    @{Row=10; A=41457; B="Sra. Alana Nogueira";   C="Vendas";     D="Viagem de negócios"; E=3; F=45091; G=4810.63},
    @{Row=11; A=66515; B="Emanuel Freitas";       C="Financeiro"; D="Outros";             E=8; F=45104; G=7348.26}
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 1).Value = $rowData.A
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
    $ws.Cells.Item($r, 7).Value = $rowData.G
}
